{"js": "// Apply the five targeted word/phrase replacements described by the diff.\n// We use Body.search() (Word JS API) with exact, case-sensitive matching so\n// each occurrence is located precisely, then replace just that hit's text\n// via Range.insertText(text, \"Replace\"). This mirrors a find & replace pass\n// done manually in the Word UI (Ctrl+H) for a handful of distinct strings.\n\nconst replacements = [\n  { find: \"social media such as messenger, or skype.\", replace: \"social media such as Messenger, or Skype.\" },\n  { find: \"sending them emails instead considering their more superior ranking.\", replace: \"sending them emails rather than considering their more superior ranking.\" },\n  { find: \"face to face rendezvous, as it all varies\", replace: \"face to face appointment, as it all varies\" },\n  { find: \"the absent crewmate or once we have received\", replace: \"the absent teammate or once we have received\" },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: \"${find}\"`);\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Apply the five targeted word/phrase replacements described by the diff\n# using Word's Find & Replace (Range.Find.Execute with Replace:=wdReplaceAll),\n# matching case so \"messenger\"/\"skype\" only hit the lower-case occurrences.\n\n$wdReplaceAll  = 2\n$wdFindStop    = 0\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text([string]$findText, [string]$replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = $wdFindStop\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute($findText, $false, $true, $false, $false, $false, $true, $wdFindStop, $false, $replaceText, $wdReplaceAll)\n    if (-not $found) {\n        throw \"Replace-Text: could not find '$findText'\"\n    }\n}\n\nReplace-Text \"social media such as messenger, or skype.\" \"social media such as Messenger, or Skype.\"\nReplace-Text \"sending them emails instead considering their more superior ranking.\" \"sending them emails rather than considering their more superior ranking.\"\nReplace-Text \"face to face rendezvous, as it all varies\" \"face to face appointment, as it all varies\"\nReplace-Text \"the absent crewmate or once we have received\" \"the absent teammate or once we have received\"\n"}
